# Bugfixed selection issue in naive_dict
# Update computed forecast-error statistics (ME, MAE, MSE, RMSE, SE, N) for rows 2-8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = -0.07723759565368317; C = 0.4832654843349234; D = 0.6010327821225036; E = 0.7752630406013843; F = 0.7951470505110364; G = 17 },
    @{ Row = 3;  B = -0.1634178345924594;  C = 0.4021168907514989; D = 0.3222171730403485; E = 0.5676417647075914; F = 0.5603403170405029; G = 17 },
    @{ Row = 4;  B = 0.05374059237319629;  C = 0.3345380109633904; D = 0.2170335242739501; E = 0.4658685697425295; F = 0.4770005923142664; G = 17 },
    @{ Row = 5;  B = -0.08511671287891998; C = 0.4691271817735259; D = 0.3831976152389424; E = 0.6190295754153774; F = 0.6320204173556079; G = 17 },
    @{ Row = 6;  B = 0.03752173280238997;  C = 0.2833056108482621; D = 0.1825500481312487; E = 0.4272587601574117; F = 0.4387066823971649; G = 17 },
    @{ Row = 7;  B = -0.01647890660520619; C = 0.3457700441494083; D = 0.1913100875655797; E = 0.4373900862680585; F = 0.4505312880675924; G = 17 },
    @{ Row = 8;  B = 0.068939684521933;    C = 0.330394643547626;  D = 0.1582419886828843; E = 0.3977964161262445; F = 0.4038346076888495; G = 17 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
}
